$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data in this sheet is stored as text (not numbers), so force the
# "@" text number format before writing any value, otherwise Excel will
# auto-convert numeric-looking strings into real numbers.
$ws.Range("A1:D6").NumberFormat = "@"

# Row 1 becomes a new header row with a single label.
$ws.Range("A1").Value = "Points"

# Existing point rows shift down by one (row N -> row N+1), and a new
# point (id 2) is appended as the last row, fixing the point ordering.
$ws.Range("A2").Value = "5"
$ws.Range("B2").Value = "-57.75229665737701"
$ws.Range("C2").Value = "-30.86397242051519"
$ws.Range("D2").Value = "0"

$ws.Range("A3").Value = "4"
$ws.Range("B3").Value = "-57.75285839699505"
$ws.Range("C3").Value = "-30.85975540560617"
$ws.Range("D3").Value = "0"

$ws.Range("A4").Value = "1"
$ws.Range("B4").Value = "-57.75526384350523"
$ws.Range("C4").Value = "-30.84392242816692"
$ws.Range("D4").Value = "0"

$ws.Range("A5").Value = "3"
$ws.Range("B5").Value = "-57.75354348580155"
$ws.Range("C5").Value = "-30.85513424867982"
$ws.Range("D5").Value = "0"

$ws.Range("A6").Value = "2"
$ws.Range("B6").Value = "-57.75428294205687"
$ws.Range("C6").Value = "-30.84993317123787"
$ws.Range("D6").Value = "0"

# Forcing text format leaves a "quote prefix" style on the touched cells;
# restore the default style now that the values are safely stored as text.
$ws.Range("A1").Style = "Normal"
$ws.Range("A2:D6").Style = "Normal"

# The old "Lines" / "Shapes" rows (6 and 7) no longer exist, and B1:D1
# must stay empty (row 1 now only has the "Points" label in column A).
$ws.Range("B1:D1").Clear()
$ws.Range("A7:D7").Clear()
